# Add a new column I to Sheet1 that concatenates each user-story row into a
# single formatted block (used for the "quick pdf" export referenced in the
# commit message), resize columns D / I, bump the row heights so the six
# wrapped lines are visible, change the zoom level, and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Widen column D slightly and add our new, very wide column I (best-fit in
# the original authoring tool so every story is readable on one "page").
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(9).ColumnWidth = 254.83333333333334

# --- Formulas in column I ------------------------------------------------
# Row 2 gets its own (non-shared) formula...
$ws.Range("I2").Formula = '=ROW()-1&". "&B2&CHAR(10)&$C$1&": "&C2&CHAR(10)&$D$1&": "&D2&CHAR(10)&$E$1&": "&E2&CHAR(10)&$F$1&": "&F2&CHAR(10)&$G$1&": "&G2'

# ...while I3:I12 share one formula (relative refs adjust row by row).
$ws.Range("I3:I12").Formula = '=ROW()-1&". "&B3&CHAR(10)&$C$1&": "&C3&CHAR(10)&$D$1&": "&D3&CHAR(10)&$E$1&": "&E3&CHAR(10)&$F$1&": "&F3&CHAR(10)&$G$1&": "&G3'

# Wrap text so the multi-line concatenation displays correctly (reuses the
# existing wrap-text cell style already used elsewhere on the sheet).
$ws.Range("I2:I12").WrapText = $true

# --- Row heights -----------------------------------------------------
$ws.Rows.Item(2).RowHeight = 99.75
$ws.Range("A3:A12").EntireRow.RowHeight = 90

# --- View / selection --------------------------------------------------
$excel.ActiveWindow.Zoom = 40
$ws.Range("I2:I11").Select()
